# StateFuncResources.xlsx edit
# Commit: "remove Gamelogic project, modify SLG building config"
#
# The sheet lists state-function effect IDs (column A) together with a
# resource/atlas id (column B). Column B only had a value for the header
# row and the "EFT_INFO" row; this change fills in the rest of column B
# with the building-config resource names ("msg_icon" for EFT_INFO's
# row, "Ssetting" for every remaining effect row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "msg_icon"
$ws.Range("B3:B15").Value = "Ssetting"
